$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cells (order matters for shared string table ordering)
$ws.Range("C14").Value = "not gabe"
$ws.Range("D6").Value = "train"
$ws.Range("A2").Value = "etse2"
$ws.Range("B1").Value = "etse3"
$ws.Range("D16").Value = $false

# Update selection to match the new active cell
$ws.Range("B2").Select()
